$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 10019
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 10019
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H29").Value = 1035.3572
$ws.Range("I29").Value = 1079.8
$ws.Range("J29").Value = 1010.6667
$ws.Range("K29").Value = 3239.4
$ws.Range("L29").Value = 3032.0001
$ws.Range("M29").Value = -2958.4
$ws.Range("N29").Value = -3594.0001
$ws.Range("H41").Value = 441.875
$ws.Range("I41").Value = 328.6
$ws.Range("J41").Value = 630.6667
$ws.Range("K41").Value = 328.6
$ws.Range("L41").Value = 630.6667
$ws.Range("M41").Value = 111.4
$ws.Range("N41").Value = -1510.6667
$ws.Range("H105").Value = 48447.5
$ws.Range("J105").Value = 48447.5
$ws.Range("L105").Value = 48447.5
$ws.Range("N105").Value = -55435.5
$ws.Range("H135").Value = 734.44446
$ws.Range("I135").Value = 395.8
$ws.Range("K135").Value = 3562.2
$ws.Range("M135").Value = -1027.2
$ws.Range("H137").Value = 838595.75
$ws.Range("I137").Value = 1430128.6
$ws.Range("K137").Value = 4290385.800000001
$ws.Range("M137").Value = -4287835.800000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4879.8
$ws.Range("J2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("N2").Value = -5226
$ws.Range("H5").Value = 372.25
$ws.Range("I5").Value = 334.83334
$ws.Range("J5").Value = 484.5
$ws.Range("K5").Value = 334.83334
$ws.Range("L5").Value = 484.5
$ws.Range("M5").Value = -222.83334
$ws.Range("N5").Value = -708.5
$ws.Range("H22").Value = 1683.2727
$ws.Range("I22").Value = 1612.8889
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1612.8889
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -1313.8889
$ws.Range("N22").Value = -2598
$ws.Range("H43").Value = 47494
$ws.Range("J43").Value = 47494
$ws.Range("L43").Value = 47494
$ws.Range("N43").Value = -48120
$ws.Range("H74").Value = 2416.9092
$ws.Range("I74").Value = 2158.6
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 2158.6
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -1284.6
$ws.Range("N74").Value = -6748
$ws.Range("H77").Value = 2416.9092
$ws.Range("I77").Value = 2158.6
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 10793
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -6425
$ws.Range("N77").Value = -33736
$ws.Range("H110").Value = 2778.1428
$ws.Range("I110").Value = 2824.5
$ws.Range("K110").Value = 2824.5
$ws.Range("M110").Value = -779.5
$ws.Range("H116").Value = 4879.8
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("N116").Value = -9588
$ws.Range("H132").Value = 9070.4
$ws.Range("I132").Value = 5905.6
$ws.Range("K132").Value = 17716.8
$ws.Range("M132").Value = -15186.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4879.8
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5228
$ws.Range("H4").Value = 372.25
$ws.Range("I4").Value = 334.83334
$ws.Range("J4").Value = 484.5
$ws.Range("K4").Value = 334.83334
$ws.Range("L4").Value = 484.5
$ws.Range("M4").Value = -219.83334
$ws.Range("N4").Value = -714.5
$ws.Range("H5").Value = 169.4
$ws.Range("I5").Value = 87.5
$ws.Range("J5").Value = 497
$ws.Range("K5").Value = 87.5
$ws.Range("L5").Value = 497
$ws.Range("M5").Value = 25.5
$ws.Range("N5").Value = -723
$ws.Range("H25").Value = 402
$ws.Range("I25").Value = 402
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 402
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -167
$ws.Range("N25").ClearContents()
$ws.Range("H105").Value = 2534.7144
$ws.Range("I105").Value = 2475.6
$ws.Range("K105").Value = 2475.6
$ws.Range("M105").Value = -728.5999999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 7833.3335
$ws.Range("J95").Value = 7833.3335
$ws.Range("L95").Value = 7833.3335
$ws.Range("N95").Value = -13325.3335
$ws.Range("H99").Value = 5892.4443
$ws.Range("I99").Value = 5717
$ws.Range("J99").Value = 6506.5
$ws.Range("K99").Value = 5717
$ws.Range("L99").Value = 6506.5
$ws.Range("M99").Value = -4219
$ws.Range("N99").Value = -9502.5
$ws.Range("H105").Value = 1951.5
$ws.Range("I105").Value = 1552.25
$ws.Range("J105").Value = 2750
$ws.Range("K105").Value = 1552.25
$ws.Range("L105").Value = 2750
$ws.Range("M105").Value = 194.75
$ws.Range("N105").Value = -6244
$ws.Range("H126").Value = 5892.4443
$ws.Range("I126").Value = 5717
$ws.Range("J126").Value = 6506.5
$ws.Range("K126").Value = 17151
$ws.Range("L126").Value = 19519.5
$ws.Range("M126").Value = -14681
$ws.Range("N126").Value = -24459.5
$ws.Range("H134").Value = 8506.076999999999
$ws.Range("I134").Value = 4016.5
$ws.Range("K134").Value = 12049.5
$ws.Range("M134").Value = -9514.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 58999.332
$ws.Range("J95").Value = 58999.332
$ws.Range("L95").Value = 58999.332
$ws.Range("N95").Value = -64491.332

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 32499
$ws.Range("J81").Value = 32499
$ws.Range("L81").Value = 32499
$ws.Range("N81").Value = -34495
$ws.Range("H84").Value = 32499
$ws.Range("J84").Value = 32499
$ws.Range("L84").Value = 97497
$ws.Range("N84").Value = -107481

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 200000
$ws.Range("I43").Value = 200000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 200000
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("M43").Value = -199851
$ws.Range("H62").Value = 4750.5
$ws.Range("I62").Value = 5334
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 5334
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -4710
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 4750.5
$ws.Range("I65").Value = 5334
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 26670
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -23550
$ws.Range("N65").Value = -21240
$ws.Range("H94").Value = 56000
$ws.Range("J94").Value = 56000
$ws.Range("L94").Value = 56000
$ws.Range("N94").Value = -57802
$ws.Range("H136").Value = 9203
$ws.Range("I136").Value = 7606
$ws.Range("J136").Value = 13195.5
$ws.Range("K136").Value = 22818
$ws.Range("L136").Value = 39586.5
$ws.Range("M136").Value = -20268
$ws.Range("N136").Value = -44686.5
